$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-31 (A,B,C,D,E,F,G,H).
# This reflects: a new sample row is inserted at the top of the series
# (pushing the ax..gz sensor readings for all prior rows down by one),
# and 10 new trailing rows are appended, extending the dataset from
# 20 rows (A1:H21) to 30 rows (A1:H31).
$data = @(
    @(0, "struggle", -0.01932663576943536, -0.1423560380935665, 0.1704595368355508, -0.0998764708638191, -0.8263479471206665, 0.256868839263916),
    @(100, "struggle", -0.002831595284598272, -0.19119707601411, 0.1513774214046346, -0.6229298114776611, -4.338823318481445, 0.1601994037628173),
    @(200, "struggle", 0.07765638828277598, -0.03652976825833296, -0.2774001099169262, -0.3721696436405182, -3.498425245285034, 0.6684392690658569),
    @(300, "struggle", 0.3138503687722347, -0.4584181170378426, -1.157653536115375, 0.1207986027002334, -1.66308057308197, 0.5348125100135803),
    @(400, "struggle", 1.636749148368837, -1.279087890471731, -4.728167861700062, -0.5253441333770752, -1.381013631820679, 0.5633704662322998),
    @(500, "struggle", 1.14729608808245, -0.2204303741455076, -3.918471813201904, 0.5025894045829773, 2.330298185348511, 0.1441642045974731),
    @(600, "struggle", -1.006189529384884, -0.457232869097163, -1.90848708152771, 0.8329147696495056, 6.17554235458374, -1.043815970420837),
    @(700, "struggle", -2.098220126969473, -1.104542408670699, 1.474607501711155, 0.0348193198442459, 5.04437255859375, 0.868803083896637),
    @(800, "struggle", -5.035840166466564, 1.3544525780848, 6.731560809271667, -0.2368630021810531, 1.756695747375489, 0.40943244099617),
    @(900, "struggle", -7.820377826690675, 2.300654053688046, 7.095898866653437, -0.0340557359158992, 2.139860868453979, 0.6606507897377014),
    @(1000, "struggle", -1.805607523236966, 2.29458946415356, 0.1832081420081089, 0.1806633919477462, 1.172861218452454, 0.4958698749542236),
    @(1100, "struggle", 1.364815167018334, 1.71517636094774, -2.800688726561401, 0.0870482996106147, -3.714518785476685, -0.6794348955154419),
    @(1200, "struggle", -2.538227507046292, -0.1809541881084444, 1.764710307121269, 0.1391245573759079, -6.393468379974365, -4.738176345825195),
    @(1300, "struggle", -2.092569351196289, 0.6709982327052525, -5.201459612165179, 1.760971784591675, -1.517847418785095, 1.018465042114258),
    @(1400, "struggle", 2.829659972872078, 1.32810748687812, -7.335261804716874, -0.0085521135479211, -4.238488674163818, -0.0195476878434419),
    @(1500, "struggle", 5.789446592330947, -0.9393511286803589, 7.550768171037943, 0.4407392740249634, -3.014314651489258, -0.6562219858169556),
    @(1600, "struggle", -0.6611437201499961, -2.4462269321084, 5.736638069152827, -0.1287398487329483, 2.51019811630249, -0.0731511116027832),
    @(1700, "struggle", -1.795438630240309, -0.6667117263589585, 0.3705058097839308, 0.2913827300071716, 6.333756446838379, 0.9952521920204164),
    @(1800, "struggle", -2.502477117947166, -0.6020157241395544, 4.68433271135602, -1.667662143707275, 4.299574851989746, 0.0829249545931816),
    @(1900, "struggle", -9.597015057291248, -0.637856924108093, 7.518706285527828, 1.54808521270752, 2.839301586151123, 1.368032693862915),
    @(2000, "struggle", -10.96845097201216, -0.236061673079225, -2.280897122408647, 0.0797179117798805, 3.36098051071167, 0.3753767013549804),
    @(2100, "struggle", 2.115621226174531, 0.9094340460641233, -5.153596741812609, 0.4167627990245819, 1.646892666816711, 0.836885392665863),
    @(2200, "struggle", -2.750307172536849, -1.468972404088292, 0.1615269269261984, 0.5007568001747131, -1.671174645423889, -0.2092213481664657),
    @(2300, "struggle", -3.695782780647288, -0.6018874943256158, 1.057165026664715, 0.0099265603348612, -3.577379703521729, -0.4952589869499206),
    @(2400, "struggle", -5.04273155757359, 2.776979684829711, -3.505371774945944, -1.279915452003479, -6.508005619049072, -0.1085812970995903),
    @(2500, "struggle", -2.32667221341812, 0.9583009992326832, -5.241142443248169, -0.9758572578430176, -1.915215253829956, 1.272432327270508),
    @(2600, "struggle", 2.433021928582876, -2.759944068534038, 3.60468020609447, 0.0806342139840126, -1.140027284622192, 0.0485637858510017),
    @(2700, "struggle", -2.125196490968979, 0.5359780830996369, 2.942476987838745, 0.8413141369819641, -0.3859141170978546, 0.5590944290161133),
    @(2800, "struggle", -3.455752406801498, 0.2274821900895648, 2.951905420848299, -0.1750128865242004, 2.383749008178711, -0.1401935666799545),
    @(2900, "struggle", -2.549259322030204, 0.480571014540536, 3.463431903294155, 0.4100432991981506, 2.355190992355347, 0.1059851199388504)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
